$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.438.46"
$ws.Range("E2").Value = "  -1.39%  "
$ws.Range("D3").Value = "1.844.03"
$ws.Range("E3").Value = "  -1.26%  "
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "265.68"
$ws.Range("E5").Value = "  -3.10%  "
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "0.5198"
$ws.Range("E7").Value = "  -1.97%  "
$ws.Range("D8").Value = "0.3270"
$ws.Range("E8").Value = "  -3.51%  "
$ws.Range("D9").Value = "0.06801"
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").Value = "18.81"
$ws.Range("E10").Value = "  -5.56%  "
$ws.Range("D11").Value = "0.7792"
$ws.Range("D12").Value = "0.07741"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").Value = "1.821.90"
$ws.Range("E13").Value = "  -1.92%  "
$ws.Range("D14").Value = "88.10"
$ws.Range("E14").Value = "  -2.28%  "
$ws.Range("D15").Value = "5.014"
$ws.Range("E15").Value = "  -2.26%  "
$ws.Range("D16").Value = "0.9991"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").Value = "13.93"
$ws.Range("E17").Value = "  -3.57%  "
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").Value = "0.000007954"
$ws.Range("E19").Value = "  -0.74%  "
$ws.Range("D20").Value = "26.452.07"
$ws.Range("E20").Value = "  -1.42%  "
$ws.Range("D21").Value = "2.073.02"
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("D22").Value = "4.626"
$ws.Range("E22").Value = "  -1.83%  "
$ws.Range("D23").Value = "9.579"
$ws.Range("E23").Value = "  -3.94%  "
$ws.Range("E24").Value = "  -1.61%  "
$ws.Range("D25").Value = "143.91"
$ws.Range("E25").Value = "  -1.15%  "
$ws.Range("D26").Value = "2.183"
$ws.Range("E26").Value = "  -8.02%  "
$ws.Range("D27").Value = "1.643"
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("D28").Value = "17.03"
$ws.Range("E28").Value = "  -1.23%  "
$ws.Range("D29").Value = "112.17"
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("D30").Value = "4.174"
$ws.Range("E30").Value = "  -3.64%  "
$ws.Range("D31").Value = "4.126"
$ws.Range("E31").Value = "  -4.47%  "
$ws.Range("D32").Value = "0.08719"
$ws.Range("E32").Value = "  -1.68%  "
$ws.Range("D33").Value = "0.04830"
$ws.Range("E33").Value = "  -1.94%  "
$ws.Range("D34").Value = "0.7230"
$ws.Range("E34").Value = "  -0.77%  "
$ws.Range("D35").Value = "1.132"
$ws.Range("E35").Value = "  -2.91%  "
$ws.Range("D36").Value = "2.843"
$ws.Range("E36").Value = "  -1.13%  "
$ws.Range("D37").Value = "3.118"
$ws.Range("E37").Value = "  -2.88%  "
$ws.Range("D38").Value = "2.234"
$ws.Range("E38").Value = "  -4.33%  "
$ws.Range("D39").Value = "0.01782"
$ws.Range("E39").Value = "  -3.69%  "
$ws.Range("D40").Value = "0.4875"
$ws.Range("E40").Value = "  -4.34%  "
$ws.Range("D41").Value = "0.9158"
$ws.Range("E41").Value = "  -2.22%  "
$ws.Range("D42").Value = "111.09"
$ws.Range("E42").Value = "  -4.64%  "
$ws.Range("D43").Value = "6.068"
$ws.Range("E43").Value = "  -0.93%  "
$ws.Range("D44").Value = "0.9998"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").Value = "7.731"
$ws.Range("E45").Value = "  -3.57%  "
$ws.Range("D46").Value = "0.4181"
$ws.Range("E46").Value = "  -5.31%  "
$ws.Range("D47").Value = "0.05927"
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("D48").Value = "9.073"
$ws.Range("E48").Value = "  -2.11%  "
$ws.Range("E49").Value = "  -6.70%  "
$ws.Range("D50").Value = "35.04"
$ws.Range("E50").Value = "  -2.92%  "
$ws.Range("D51").Value = "0.8862"
$ws.Range("E51").Value = "  +0.77%  "
